$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up the existing R-column cells (row 4-7) so they use the same
# styles as their left neighbours instead of the special one-off styles
# that are being retired from the workbook's style table.
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)

$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial(-4122)

$ws.Range("Q6").Copy()
$ws.Range("R6").PasteSpecial(-4122)

$ws.Range("Q7").Copy()
$ws.Range("R7").PasteSpecial(-4122)

# --- Add the new 2022 column (S) with the same formatting as column R.
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)
$ws.Range("S4").Value = 2022

$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)
$ws.Range("S5").Value = 49.7

$ws.Range("R6").Copy()
$ws.Range("S6").PasteSpecial(-4122)
$ws.Range("S6").Value = 34.9

$ws.Range("R7").Copy()
$ws.Range("S7").PasteSpecial(-4122)
$ws.Range("S7").Value = 21

# --- Update the active selection to match the saved view state.
$null = $ws.Range("R12").Select()
